$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4357.0835
$ws.Range("I11").Value = 4357.0835
$ws.Range("K11").Value = 4357.0835
$ws.Range("M11").Value = -4217.0835
$ws.Range("H121").Value = 975
$ws.Range("J121").Value = 975
$ws.Range("L121").Value = 2925
$ws.Range("N121").Value = -6419
$ws.Range("H129").Value = 1210.9395
$ws.Range("I129").Value = 259.5
$ws.Range("J129").Value = 1422.3704
$ws.Range("K129").Value = 778.5
$ws.Range("L129").Value = 4267.1112
$ws.Range("M129").Value = 4221.5
$ws.Range("N129").Value = -14267.1112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4537.6753
$ws.Range("I32").Value = 4401.3184
$ws.Range("J32").Value = 5355.8184
$ws.Range("K32").Value = 4401.3184
$ws.Range("L32").Value = 5355.8184
$ws.Range("M32").Value = -4114.3184
$ws.Range("N32").Value = -5929.8184
$ws.Range("H61").Value = 2848.261
$ws.Range("I61").Value = 2350.5557
$ws.Range("J61").Value = 4640
$ws.Range("K61").Value = 2350.5557
$ws.Range("L61").Value = 4640
$ws.Range("M61").Value = -2138.5557
$ws.Range("N61").Value = -5064
$ws.Range("H74").Value = 8450.666999999999
$ws.Range("I74").Value = 1182.5
$ws.Range("J74").Value = 16757.143
$ws.Range("K74").Value = 1182.5
$ws.Range("L74").Value = 16757.143
$ws.Range("M74").Value = -308.5
$ws.Range("N74").Value = -18505.143
$ws.Range("H77").Value = 8450.666999999999
$ws.Range("I77").Value = 1182.5
$ws.Range("J77").Value = 16757.143
$ws.Range("K77").Value = 5912.5
$ws.Range("L77").Value = 83785.715
$ws.Range("M77").Value = -1544.5
$ws.Range("N77").Value = -92521.715
$ws.Range("H122").Value = 1156.3158
$ws.Range("I122").Value = 1051.1765
$ws.Range("K122").Value = 3153.5295
$ws.Range("M122").Value = -703.5295000000001
$ws.Range("H132").Value = 1462.1086
$ws.Range("I132").Value = 1329.6923
$ws.Range("J132").Value = 2199.8572
$ws.Range("K132").Value = 3989.0769
$ws.Range("L132").Value = 6599.571599999999
$ws.Range("M132").Value = -1459.0769
$ws.Range("N132").Value = -11659.5716
$ws.Range("H136").Value = 2848.261
$ws.Range("I136").Value = 2350.5557
$ws.Range("J136").Value = 4640
$ws.Range("K136").Value = 7051.6671
$ws.Range("L136").Value = 13920
$ws.Range("M136").Value = -4501.6671
$ws.Range("N136").Value = -19020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1600
$ws.Range("I29").Value = 1600
$ws.Range("K29").Value = 1600
$ws.Range("M29").Value = -1311
$ws.Range("H134").Value = 45428.957
$ws.Range("I134").Value = 51952.6
$ws.Range("K134").Value = 155857.8
$ws.Range("M134").Value = -153322.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1214.1666
$ws.Range("I31").Value = 1192
$ws.Range("J31").Value = 1325
$ws.Range("K31").Value = 1192
$ws.Range("L31").Value = 1325
$ws.Range("M31").Value = -897
$ws.Range("N31").Value = -1915
$ws.Range("H34").Value = 1214.1666
$ws.Range("I34").Value = 1192
$ws.Range("J34").Value = 1325
$ws.Range("K34").Value = 1192
$ws.Range("L34").Value = 1325
$ws.Range("M34").Value = -990
$ws.Range("N34").Value = -1729
$ws.Range("H112").Value = 13801
$ws.Range("J112").Value = 13801
$ws.Range("L112").Value = 13801
$ws.Range("N112").Value = -16755
$ws.Range("H132").Value = 1962.6471
$ws.Range("I132").Value = 1846.1143
$ws.Range("J132").Value = 2217.5625
$ws.Range("K132").Value = 5538.3429
$ws.Range("L132").Value = 6652.6875
$ws.Range("M132").Value = -3008.3429
$ws.Range("N132").Value = -11712.6875
$ws.Range("H134").Value = 2043.7317
$ws.Range("I134").Value = 2292.6562
$ws.Range("J134").Value = 1158.6666
$ws.Range("K134").Value = 6877.9686
$ws.Range("L134").Value = 3475.9998
$ws.Range("M134").Value = -4342.9686
$ws.Range("N134").Value = -8545.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43550
$ws.Range("I12").Value = 112
$ws.Range("J12").Value = 71474.42999999999
$ws.Range("K12").Value = 336
$ws.Range("L12").Value = 214423.29
$ws.Range("M12").Value = -163
$ws.Range("N12").Value = -214769.29
$ws.Range("H59").Value = 2890
$ws.Range("I59").Value = 2800
$ws.Range("J59").Value = 2980
$ws.Range("K59").Value = 8400
$ws.Range("L59").Value = 8940
$ws.Range("M59").Value = -7860
$ws.Range("N59").Value = -10020
$ws.Range("H80").Value = 9054692
$ws.Range("J80").Value = 1594.6
$ws.Range("L80").Value = 4783.799999999999
$ws.Range("N80").Value = -6655.799999999999
$ws.Range("H83").Value = 9054692
$ws.Range("J83").Value = 1594.6
$ws.Range("L83").Value = 14351.4
$ws.Range("N83").Value = -23711.4
$ws.Range("H94").Value = 2308
$ws.Range("J94").Value = 2800
$ws.Range("L94").Value = 8400
$ws.Range("N94").Value = -9752
$ws.Range("H96").Value = 70707200
$ws.Range("J96").Value = 70707200
$ws.Range("L96").Value = 212121600
$ws.Range("N96").Value = -212125718
$ws.Range("H98").Value = 490.08334
$ws.Range("I98").Value = 413
$ws.Range("K98").Value = 1239
$ws.Range("M98").Value = 259
$ws.Range("H107").Value = 61020.816
$ws.Range("J107").Value = 38910.92
$ws.Range("L107").Value = 116732.76
$ws.Range("N107").Value = -120572.76
$ws.Range("H110").Value = 3514.2856
$ws.Range("I110").Value = 3150
$ws.Range("J110").Value = 3660
$ws.Range("K110").Value = 9450
$ws.Range("L110").Value = 10980
$ws.Range("M110").Value = -5360
$ws.Range("N110").Value = -19160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 1292.5927
$ws.Range("I24").Value = 2725
$ws.Range("J24").Value = 1043.4783
$ws.Range("K24").Value = 2725
$ws.Range("L24").Value = 1043.4783
$ws.Range("M24").Value = -2552
$ws.Range("N24").Value = -1389.4783

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1781.0714
$ws.Range("I68").Value = 1693.5
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1693.5
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -944.5
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1781.0714
$ws.Range("I71").Value = 1693.5
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 8467.5
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -4723.5
$ws.Range("N71").Value = -17488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("N18").Value = -3346
$ws.Range("H22").Value = 2950
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1707
$ws.Range("H113").Value = 661.3333
$ws.Range("I113").Value = 866.6667
$ws.Range("J113").Value = 558.6667
$ws.Range("K113").Value = 2600.0001
$ws.Range("L113").Value = 1676.0001
$ws.Range("M113").Value = -430.0001000000002
$ws.Range("N113").Value = -6016.0001
$ws.Range("H122").Value = 1937.2821
$ws.Range("I122").Value = 1299
$ws.Range("J122").Value = 3077.0715
$ws.Range("K122").Value = 3897
$ws.Range("L122").Value = 9231.2145
$ws.Range("M122").Value = -1447
$ws.Range("N122").Value = -14131.2145
